# Update row 8 (ano=2025) metrics in metricas_recorrencia_anual
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1125
$ws.Range("D8").Value = 185
$ws.Range("E8").Value = 940
$ws.Range("F8").Value = 7.588187038556193
$ws.Range("G8").Value = 83.55555555555556
$ws.Range("H8").Value = 16.44444444444445
